$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178327441215515
$ws.Range("B1").Value = 2.418038368225098
$ws.Range("D1").Value = 2.331649780273438
$ws.Range("E1").Value = 1.197895765304565
